$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.732.02"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.891.56"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4881"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06672"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "1.886.59"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07251"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "89.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.013"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6645"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").Value = "30.681.99"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007885"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "2.130.76"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.742"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.065"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.823"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.267"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09032"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05208"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7324"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.086"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.696"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01822"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9232"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.061"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4446"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.725"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.310"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4226"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05837"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.649"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
